$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 37
$ws.Cells.Item($row, 1).Value = "05/01/2026 02:24:12"
$ws.Cells.Item($row, 2).Value = "05/01 02:00"
$ws.Cells.Item($row, 3).Value = "Metrópoles"
$ws.Cells.Item($row, 4).Value = "Entenda PL que corta benefícios fiscais e aumenta tributação de bets"
$ws.Cells.Item($row, 5).Value = "https://www.metropoles.com/brasil/economia-br/entenda-pl-que-corta-beneficios-fiscais-e-aumenta-tributacao-de-bets"
$ws.Cells.Item($row, 6).Value = "imposto"
$ws.Cells.Item($row, 7).Value = "da busca reforçar a arrecadação e fechar as contas de 2026 ao reduzir incentivos e elevar &lt;b&gt;imposto&lt;/b&gt;s sobre bets e fintechs"
